$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("output_and_resourceuse")
$ws2 = $wb.Worksheets.Item("optimal_coverage")

# --- Sheet1 (output_and_resourceuse): update row 8 with new HR-constraint scenario ---
$ws1.Cells.Item(8, 1).Value = "CET ($164.7) + Demand constraint + Drug budget + HR constraint"
$ws1.Cells.Item(8, 2).Value = 89
$ws1.Cells.Item(8, 3).Value = 56
$ws1.Cells.Item(8, 4).Value = 64412836.83
$ws1.Cells.Item(8, 5).Value = 152.79
$ws1.Cells.Item(8, 6).Value = 0.9892513549003035
$ws1.Cells.Item(8, 7).Value = 1
$ws1.Cells.Item(8, 8).Value = 0.43
$ws1.Cells.Item(8, 9).Value = 0.85
$ws1.Cells.Item(8, 10).Value = 0.21

# --- Sheet2 (optimal_coverage): rename header for last scenario column (J1) ---
$ws2.Cells.Item(1, 10).Value = "CET ($164.7) + Demand constraint + Drug budget + HR constraint"

# --- Sheet2: fix column order so A=Program, B=Intervention code (matches header) ---
for ($r = 2; $r -le 142; $r++) {
    $codeCell = $ws2.Cells.Item($r, 1)
    $progCell = $ws2.Cells.Item($r, 2)
    $code = $codeCell.Value()
    $prog = $progCell.Value()

    $codeCell.Value = $prog
    $progCell.NumberFormat = "@"
    $progCell.Value = $code
    $progCell.ClearFormats()
}

# --- Sheet2: refresh the "Donor constraint" column (J) values with new HR-constraint results ---
$jUpdates = @{
    2 = 0.4699999999999921
    5 = 0.7100000000004439
    7 = 0
    8 = 0.9999999999993293
    10 = 0.6999999999999998
    11 = 1.000000000000037
    12 = 0
    13 = 0
    18 = 0.4
    19 = 0.6
    20 = 0
    21 = 0.6000000000000001
    22 = 0.25
    23 = 0
    24 = 0
    26 = 0.8999999999999999
    30 = 0.7
    31 = 0
    32 = 0.9999999999996588
    34 = 0.6000000000001343
    35 = 0.6000000000002146
    36 = 0.5999999999999478
    38 = 0.95
    39 = 0
    40 = 0
    41 = 0
    42 = 0.8000000000000136
    43 = 0.9999999999994368
    45 = 0.9499999999999997
    46 = 0.95
    48 = 0
    50 = 0.6000000000001128
    51 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0.5000000000003847
    59 = 0.8000000000006154
    60 = 0.05000000000000001
    63 = 0.9500000000000001
    65 = 0.9499999999993444
    66 = 0.9499999999994824
    67 = 0.95
    69 = 0
    70 = 0
    71 = 0.9499999999997762
    81 = 0.6000000000002113
    86 = 1
    91 = 0.9000000000006638
    95 = 0
    96 = 0
    100 = 0.9999999999994607
    101 = 0.2005600796501301
    104 = 0.8999999999997979
    107 = 0
    109 = 0
    114 = 0.9500000000000001
    122 = 0.9
    123 = 0
    124 = 0.92
    127 = 0.7999999999999998
    128 = 0
    129 = 0.9800000000000001
    130 = 0
    131 = 1
    132 = 0.97
    134 = 0.9800000000000001
    135 = 0
    137 = 0.8000000000001412
    138 = 0.6000000000002502
    139 = 0.6000000000000001
    140 = 0.59
}
foreach ($r in $jUpdates.Keys) {
    $ws2.Cells.Item($r, 10).Value = $jUpdates[$r]
}
